$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet Property1 -> DataNode
$ws.Name = "DataNode"
